$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 41
$ws.Cells.Item(41, 8).Value = 1159.7894
$ws.Cells.Item(41, 9).Value = 1400.5
$ws.Cells.Item(41, 10).Value = 892.3333
$ws.Cells.Item(41, 11).Value = 1400.5
$ws.Cells.Item(41, 12).Value = 892.3333
$ws.Cells.Item(41, 13).Value = -960.5
$ws.Cells.Item(41, 14).Value = -1772.3333

# ALC row 64
$ws.Cells.Item(64, 8).Value = 6172.607
$ws.Cells.Item(64, 10).Value = 7585.2144
$ws.Cells.Item(64, 12).Value = 7585.2144
$ws.Cells.Item(64, 14).Value = -8081.2144

# ALC row 67
$ws.Cells.Item(67, 8).Value = 6172.607
$ws.Cells.Item(67, 10).Value = 7585.2144
$ws.Cells.Item(67, 12).Value = 7585.2144
$ws.Cells.Item(67, 14).Value = -9301.214400000001

# ALC row 74
$ws.Cells.Item(74, 8).Value = 5008.476
$ws.Cells.Item(74, 9).Value = 4811.6
$ws.Cells.Item(74, 11).Value = 4811.6
$ws.Cells.Item(74, 13).Value = -3875.6

# ALC row 77
$ws.Cells.Item(77, 8).Value = 5008.476
$ws.Cells.Item(77, 9).Value = 4811.6
$ws.Cells.Item(77, 11).Value = 24058
$ws.Cells.Item(77, 13).Value = -19378

# ALC row 133
$ws.Cells.Item(133, 8).Value = 76107.164
$ws.Cells.Item(133, 10).Value = 76107.164
$ws.Cells.Item(133, 12).Value = 76107.164
$ws.Cells.Item(133, 14).Value = -86227.164

# ALC row 134
$ws.Cells.Item(134, 8).Value = 39996.273
$ws.Cells.Item(134, 10).Value = 39996.273
$ws.Cells.Item(134, 12).Value = 39996.273
$ws.Cells.Item(134, 14).Value = -50136.273

# ALC row 136
$ws.Cells.Item(136, 8).Value = 57579.5
$ws.Cells.Item(136, 10).Value = 57579.5
$ws.Cells.Item(136, 12).Value = 57579.5
$ws.Cells.Item(136, 14).Value = -67779.5

# ALC row 137
$ws.Cells.Item(137, 8).Value = 502434.53
$ws.Cells.Item(137, 9).Value = 1480.35
$ws.Cells.Item(137, 11).Value = 4441.049999999999
$ws.Cells.Item(137, 13).Value = -1891.049999999999

# ALC row 139
$ws.Cells.Item(139, 8).Value = 73843
$ws.Cells.Item(139, 10).Value = 73843
$ws.Cells.Item(139, 12).Value = 73843
$ws.Cells.Item(139, 14).Value = -84123

$ws = $wb.Worksheets.Item("ARM")
# ARM row 32
$ws.Cells.Item(32, 8).Value = 13483.571
$ws.Cells.Item(32, 9).Value = 6007.4287
$ws.Cells.Item(32, 11).Value = 6007.4287
$ws.Cells.Item(32, 13).Value = -5720.4287

# ARM row 33
$ws.Cells.Item(33, 8).Value = 18750
$ws.Cells.Item(33, 9).Value = 19500
$ws.Cells.Item(33, 10).Value = 18000
$ws.Cells.Item(33, 11).Value = 19500
$ws.Cells.Item(33, 12).Value = 18000
$ws.Cells.Item(33, 13).Value = -19171
$ws.Cells.Item(33, 14).Value = -18658

# ARM row 74
$ws.Cells.Item(74, 8).Value = 26763.35
$ws.Cells.Item(74, 9).Value = 31401.121
$ws.Cells.Item(74, 10).Value = 4899.5713
$ws.Cells.Item(74, 11).Value = 31401.121
$ws.Cells.Item(74, 12).Value = 4899.5713
$ws.Cells.Item(74, 13).Value = -30527.121
$ws.Cells.Item(74, 14).Value = -6647.5713

# ARM row 77
$ws.Cells.Item(77, 8).Value = 26763.35
$ws.Cells.Item(77, 9).Value = 31401.121
$ws.Cells.Item(77, 10).Value = 4899.5713
$ws.Cells.Item(77, 11).Value = 157005.605
$ws.Cells.Item(77, 12).Value = 24497.8565
$ws.Cells.Item(77, 13).Value = -152637.605
$ws.Cells.Item(77, 14).Value = -33233.85649999999

# ARM row 102
$ws.Cells.Item(102, 8).Value = 94643.914
$ws.Cells.Item(102, 9).Value = 126492.5
$ws.Cells.Item(102, 10).Value = 30946.75
$ws.Cells.Item(102, 11).Value = 126492.5
$ws.Cells.Item(102, 12).Value = 30946.75
$ws.Cells.Item(102, 13).Value = -124870.5
$ws.Cells.Item(102, 14).Value = -34190.75

# ARM row 110
$ws.Cells.Item(110, 8).Value = 2376.8333
$ws.Cells.Item(110, 9).Value = 1999
$ws.Cells.Item(110, 11).Value = 1999
$ws.Cells.Item(110, 13).Value = 46

# ARM row 118
$ws.Cells.Item(118, 8).Value = 54497.5
$ws.Cells.Item(118, 10).Value = 54497.5
$ws.Cells.Item(118, 12).Value = 54497.5
$ws.Cells.Item(118, 14).Value = -57811.5

# ARM row 122
$ws.Cells.Item(122, 8).Value = 3370.5
$ws.Cells.Item(122, 9).Value = 3189.4443
$ws.Cells.Item(122, 11).Value = 9568.332900000001
$ws.Cells.Item(122, 13).Value = -7118.332900000001

# ARM row 127
$ws.Cells.Item(127, 8).Value = 94996.664
$ws.Cells.Item(127, 10).Value = 94996.664
$ws.Cells.Item(127, 12).Value = 94996.664
$ws.Cells.Item(127, 14).Value = -104916.664

$ws = $wb.Worksheets.Item("BSM")
# BSM row 52
$ws.Cells.Item(52, 8).Value = 99990
$ws.Cells.Item(52, 10).Value = 99990
$ws.Cells.Item(52, 12).Value = 99990
$ws.Cells.Item(52, 14).Value = -100516

# BSM row 80
$ws.Cells.Item(80, 8).Value = 2447.625
$ws.Cells.Item(80, 10).Value = 2558.2
$ws.Cells.Item(80, 12).Value = 2558.2
$ws.Cells.Item(80, 14).Value = -4554.2

# BSM row 83
$ws.Cells.Item(83, 8).Value = 2447.625
$ws.Cells.Item(83, 10).Value = 2558.2
$ws.Cells.Item(83, 12).Value = 12791
$ws.Cells.Item(83, 14).Value = -22775

# BSM row 108
$ws.Cells.Item(108, 8).Value = 99995
$ws.Cells.Item(108, 10).Value = 99995
$ws.Cells.Item(108, 12).Value = 99995
$ws.Cells.Item(108, 14).Value = -107675

# BSM row 115
$ws.Cells.Item(115, 8).Value = 76996.836
$ws.Cells.Item(115, 10).Value = 79996
$ws.Cells.Item(115, 12).Value = 79996
$ws.Cells.Item(115, 14).Value = -83130

# BSM row 119
$ws.Cells.Item(119, 8).Value = 99392
$ws.Cells.Item(119, 10).Value = 99392
$ws.Cells.Item(119, 12).Value = 99392
$ws.Cells.Item(119, 14).Value = -109068

# BSM row 121
$ws.Cells.Item(121, 8).Value = 99990
$ws.Cells.Item(121, 10).Value = 99990
$ws.Cells.Item(121, 12).Value = 99990
$ws.Cells.Item(121, 14).Value = -103484

# BSM row 127
$ws.Cells.Item(127, 8).Value = 57484
$ws.Cells.Item(127, 10).Value = 57484
$ws.Cells.Item(127, 12).Value = 57484
$ws.Cells.Item(127, 14).Value = -67404

# BSM row 132
$ws.Cells.Item(132, 8).Value = 30292.195
$ws.Cells.Item(132, 10).Value = 30292.195
$ws.Cells.Item(132, 12).Value = 30292.195
$ws.Cells.Item(132, 14).Value = -40412.195

# BSM row 134
$ws.Cells.Item(134, 8).Value = 3204.1282
$ws.Cells.Item(134, 9).Value = 2079.5806
$ws.Cells.Item(134, 11).Value = 6238.7418
$ws.Cells.Item(134, 13).Value = -3703.7418

# BSM row 138
$ws.Cells.Item(138, 8).Value = 74652.25
$ws.Cells.Item(138, 10).Value = 74652.25
$ws.Cells.Item(138, 12).Value = 74652.25
$ws.Cells.Item(138, 14).Value = -84932.25

# BSM row 140
$ws.Cells.Item(140, 8).Value = 42330.668
$ws.Cells.Item(140, 10).Value = 42330.668
$ws.Cells.Item(140, 12).Value = 42330.668
$ws.Cells.Item(140, 14).Value = -52690.668

$ws = $wb.Worksheets.Item("CRP")
# CRP row 3
$ws.Cells.Item(3, 8).Value = 0
$ws.Cells.Item(3, 9).Value = 0
$ws.Cells.Item(3, 10).Value = 0
$ws.Cells.Item(3, 11).Value = 0
$ws.Cells.Item(3, 12).Value = 0
$ws.Cells.Item(3, 13).ClearContents()
$ws.Cells.Item(3, 14).ClearContents()

# CRP row 62
$ws.Cells.Item(62, 8).Value = 3518.4
$ws.Cells.Item(62, 9).Value = 3864
$ws.Cells.Item(62, 10).Value = 3000
$ws.Cells.Item(62, 11).Value = 3864
$ws.Cells.Item(62, 12).Value = 3000
$ws.Cells.Item(62, 13).Value = -3240
$ws.Cells.Item(62, 14).Value = -4248

# CRP row 65
$ws.Cells.Item(65, 8).Value = 3518.4
$ws.Cells.Item(65, 9).Value = 3864
$ws.Cells.Item(65, 10).Value = 3000
$ws.Cells.Item(65, 11).Value = 19320
$ws.Cells.Item(65, 12).Value = 15000
$ws.Cells.Item(65, 13).Value = -16200
$ws.Cells.Item(65, 14).Value = -21240

# CRP row 114
$ws.Cells.Item(114, 8).Value = 39971.5
$ws.Cells.Item(114, 10).Value = 39971.5
$ws.Cells.Item(114, 12).Value = 39971.5
$ws.Cells.Item(114, 14).Value = -48649.5

# CRP row 118
$ws.Cells.Item(118, 8).Value = 64997.5
$ws.Cells.Item(118, 10).Value = 64997.5
$ws.Cells.Item(118, 12).Value = 64997.5
$ws.Cells.Item(118, 14).Value = -68311.5

# CRP row 138
$ws.Cells.Item(138, 8).Value = 99996
$ws.Cells.Item(138, 10).Value = 99996
$ws.Cells.Item(138, 12).Value = 99996
$ws.Cells.Item(138, 14).Value = -110276

$ws = $wb.Worksheets.Item("GSM")
# GSM row 5
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 13).ClearContents()

# GSM row 93
$ws.Cells.Item(93, 8).Value = 14117.8
$ws.Cells.Item(93, 10).Value = 14117.8
$ws.Cells.Item(93, 12).Value = 14117.8
$ws.Cells.Item(93, 14).Value = -17861.8

# GSM row 109
$ws.Cells.Item(109, 8).Value = 23999.363
$ws.Cells.Item(109, 10).Value = 23999.363
$ws.Cells.Item(109, 12).Value = 23999.363
$ws.Cells.Item(109, 14).Value = -26079.363

# GSM row 114
$ws.Cells.Item(114, 8).Value = 79489.25
$ws.Cells.Item(114, 10).Value = 79489.25
$ws.Cells.Item(114, 12).Value = 79489.25
$ws.Cells.Item(114, 14).Value = -88167.25

# GSM row 116
$ws.Cells.Item(116, 8).Value = 57780.375
$ws.Cells.Item(116, 10).Value = 58892
$ws.Cells.Item(116, 12).Value = 58892
$ws.Cells.Item(116, 14).Value = -68070

# GSM row 119
$ws.Cells.Item(119, 8).Value = 59996.8
$ws.Cells.Item(119, 10).Value = 59996.8
$ws.Cells.Item(119, 12).Value = 59996.8
$ws.Cells.Item(119, 14).Value = -69672.8

# GSM row 126
$ws.Cells.Item(126, 8).Value = 5510
$ws.Cells.Item(126, 9).Value = 4000
$ws.Cells.Item(126, 10).Value = 6516.6665
$ws.Cells.Item(126, 11).Value = 12000
$ws.Cells.Item(126, 12).Value = 19549.9995
$ws.Cells.Item(126, 13).Value = -9530
$ws.Cells.Item(126, 14).Value = -24489.9995

# GSM row 135
$ws.Cells.Item(135, 8).Value = 45436.25
$ws.Cells.Item(135, 10).Value = 45436.25
$ws.Cells.Item(135, 12).Value = 45436.25
$ws.Cells.Item(135, 14).Value = -55576.25

# GSM row 140
$ws.Cells.Item(140, 8).Value = 95206
$ws.Cells.Item(140, 10).Value = 95633.82000000001
$ws.Cells.Item(140, 12).Value = 95633.82000000001
$ws.Cells.Item(140, 14).Value = -105993.82

$ws = $wb.Worksheets.Item("LTW")
# LTW row 117
$ws.Cells.Item(117, 8).Value = 89096
$ws.Cells.Item(117, 10).Value = 89096
$ws.Cells.Item(117, 12).Value = 89096
$ws.Cells.Item(117, 14).Value = -98274

# LTW row 118
$ws.Cells.Item(118, 8).Value = 98401.60000000001
$ws.Cells.Item(118, 10).Value = 98401.60000000001
$ws.Cells.Item(118, 12).Value = 98401.60000000001
$ws.Cells.Item(118, 14).Value = -101715.6

# LTW row 127
$ws.Cells.Item(127, 8).Value = 250000
$ws.Cells.Item(127, 10).Value = 250000
$ws.Cells.Item(127, 12).Value = 250000
$ws.Cells.Item(127, 14).Value = -259920

$ws = $wb.Worksheets.Item("WVR")
# WVR row 127
$ws.Cells.Item(127, 8).Value = 60390
$ws.Cells.Item(127, 9).Value = 60390
$ws.Cells.Item(127, 11).Value = 60390
$ws.Cells.Item(127, 13).Value = -55430
